# Updated MCH102 to MCH251
# Add a new data row (row 2) under the existing header row with the
# archival-description metadata for MCH128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MCH128"
$ws.Range("C2").Value = "CORRESPONDENCE, LET MY PEOPLE GO, THE PEOPLE SHALL GOVERN"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

$row2 = $ws.Range("A2:H2")
$row2.Font.Name = "Calibri"
$row2.Font.ThemeColor = 1

$ws.Range("B2").Clear()

$ws.Range("A2:H2").Select()
$excel.ActiveWindow.FreezePanes = $true
